$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 15: "ADC Grenze" (threshold between adjacent ADC steps) ---
$ws.Range("B15").Value = "ADC Grenze"

# C15 and J15 are "edge" cells with their own one-off formula, so set them
# individually first...
$ws.Range("C15").Formula = "=((C13+D13)/2)"
$ws.Range("J15").Formula = "=((J13+C10)/2)"

# ...then fill D15:I15 in one shot so the engine records them as a single
# shared-formula block (matches Excel's own fill-right behaviour).
$ws.Range("D15:I15").Formula = "=((D13+E13)/2)"

# --- Narrower, uniform column widths for C:J (was bestFit 12.23 / 13.23) ---
$ws.Range("C1:J1").ColumnWidth = 7.8

# --- Selection moved to M10 ---
[void]$ws.Range("M10").Select()
